# Add a new review row (row 6) to Sheet1, mirroring the structure/formatting
# of the existing rows, with two new hyperlinked e-mail addresses and a new
# review comment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 5 (values + formatting) into a new row 6, so the new row
# starts out with the same styles/number formats as the rest of the table.
$ws.Rows("5:5").Copy()
$ws.Rows("6:6").Insert()

# Register the two new hyperlinks (this also happens to tweak the cell
# formatting of C6/D6, which we restore below).
$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:cybworking@gmail.com", $null, $null, "cybworking@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:fadeaway12222@gmail.com", $null, $null, "fadeaway12222@gmail.com")

# Fill in the new row's data.
$ws.Range("C6").Value = "cybworking@gmail.com"
$ws.Range("D6").Value = "fadeaway12222@gmail.com"
$ws.Range("F6").Value = "I think this game is excellent for kids and adults. Very fun to play it along the day…"

# Re-apply row 5's formatting across the whole new row so that every cell
# (including C6/D6, whose format was altered by adding the hyperlinks)
# keeps the same look as the rest of the table.
$ws.Range("A5:F5").Copy()
$ws.Range("A6").PasteSpecial(-4122)

# Clean up the now-unused "Hyperlink" cell style that got created as a
# side-effect of Hyperlinks.Add, since it is no longer referenced anywhere.
$wb.Styles.Item("Hyperlink").Delete()

# Match the workbook's saved selection / active cell on the new last cell.
$ws.Range("F6").Select()
